# dataMen.xlsx: the J:L and M:O column blocks (sib1_sex1_* / sib1_sex0_*
# stats) were swapped - both the header labels in row 1 and the data
# values in every data row - so that the sib1_sex1_* block now sits in
# J:L and the sib1_sex0_* block now sits in M:O.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where both the J:L block and the M:O block already contain data -
# a straight value swap between the two blocks.
$fullSwapRows = @(1,2,3,4,5,6,7,8,9,10,11,12,13,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,52,53,54,56,57,58,59,63,64,72,73,74,77,85,86,87,88,90,91,92,93)

foreach ($r in $fullSwapRows) {
    $jkl = $ws.Range("J$r`:L$r").Value()
    $mno = $ws.Range("M$r`:O$r").Value()
    $ws.Range("J$r`:L$r").Value = $mno
    $ws.Range("M$r`:O$r").Value = $jkl
}

# Rows where only the J:L block has data and the M:O block is blank -
# move the values over to M:O and clear out J:L entirely.
$moveOnlyRows = @(14,30,31,49,51,61,82)

foreach ($r in $moveOnlyRows) {
    $jkl = $ws.Range("J$r`:L$r").Value()
    $ws.Range("M$r`:O$r").Value = $jkl
    $ws.Range("J$r`:L$r").Clear()
}
